$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "decision_date" column header (new shared string)
$ws.Range("D1").Value = "decision_date"

# Fill D2:D210 with the decision date formula (copom meeting date minus 8 days)
$ws.Range("D2:D210").Formula = "=A2-8"

# Prepare the new (currently empty) I5:J5 block: copy the date number-format
# from an existing date cell so the cells pick up the same style (s=1) as
# column A/D without introducing a brand-new number format definition.
$ws.Range("A2").Copy()
$ws.Range("I5:J5").PasteSpecial(-4122)

# Leave the selection where the author left it after the edit
$ws.Range("I5:J5").Select()
